$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 487, shifting existing rows 487-597 down to 488-598.
$ws.Rows.Item(487).Insert()

# Populate the newly inserted row 487 with the new record's data.
$ws.Range("A487").Value = 8
$ws.Range("B487").Value = "Terminal La Palmera de La Serena"
$ws.Range("C487").Value = "Coquimbo"
$ws.Range("D487").Value = 45173
$ws.Range("D487").NumberFormat = $ws.Range("D488").NumberFormat
$ws.Range("E487").Value = 4
$ws.Range("F487").Value = 100114013
$ws.Range("G487").Value = "Zanahoria"
$ws.Range("H487").Value = "Sin especificar"
$ws.Range("I487").Value = "Primera"
$ws.Range("J487").Value = 520
$ws.Range("K487").Value = 5500
$ws.Range("L487").Value = 6000
$ws.Range("M487").Value = 5750
$ws.Range("N487").Value = "$/saco 20 kilos"
$ws.Range("O487").Value = "Provincia del Elquí"
$ws.Range("P487").Value = 288
$ws.Range("Q487").Value = 20
$ws.Range("R487").Value = "Hortaliza"
